# "Updating the model for Horeco"
# Rolls the whole Consumption_Actual_Historical sheet forward by one day:
#   - every 15-min Timestamp in column A moves from day D to day D+1
#   - column B (Actual Consumption) is replaced with the newly observed
#     values for the corresponding quarter-hour of the new day
#   - column D (Lookup = "dd.mm.yyyy" + quarter-of-day index) is
#     regenerated from the shifted date in column A and the existing
#     quarter counter already stored in column C
#   - column C (Quarter) is untouched

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly observed "Actual Consumption (MW)" readings for the day that is
# being appended (row 2 = first quarter-hour .. row 193 = last quarter-hour
# two days later), in row order.
$newB = @(6143,6056,6069,6008,6033,5985,6030,5993,6012,5999,6039,6017,6038,6008,6100,6105,6208,6205,6297,6418,6605,6745,6845,7008,7211,7400,7481,7647,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,6261,6206,6209,6166,6104,6065,6056,6038,5994,5969,5990,5980,5994,5976,5992,6006,6040,6061,6095,6130,6167,6202,6241,6294,6289,6472,6551,6735,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$startRow = 2
$endRow = 193

for ($r = $startRow; $r -le $endRow; $r++) {
    $idx = $r - $startRow

    # Column A: timestamp shifts forward by exactly one day, keeping the
    # same time-of-day fraction.
    $curSerial = $ws.Cells.Item($r, 1).Value2()
    $ws.Cells.Item($r, 1).Value2 = $curSerial + 1

    # Column B: actual consumption value fetched for the new day.
    $ws.Cells.Item($r, 2).Value = $newB[$idx]

    # Column D: rebuild the "dd.mm.yyyy" + quarter-index lookup label from
    # the (now shifted) date in column A and the quarter counter in C.
    $newDate = $ws.Cells.Item($r, 1).Value()
    $quarter = $ws.Cells.Item($r, 3).Value2()
    $day = "{0:D2}" -f $newDate.Day
    $month = "{0:D2}" -f $newDate.Month
    $year = $newDate.Year
    $ws.Cells.Item($r, 4).Value = "$day.$month.$year$quarter"
}
